# Weekly update: insert a new price record at the top of the Mango data
# table (row 151), shifting all existing rows down by one. This mirrors
# how the source system prepends the latest week's observation while
# keeping historical rows intact (the oldest row simply moves to the
# new last row of the range).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 151; rows 151:200 shift to 152:201.
$ws.Rows(151).Insert()

# Populate the newly inserted row with this week's observation.
$ws.Cells.Item(151, 1).Value = 7
$ws.Cells.Item(151, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(151, 3).Value = "Ñuble"
$ws.Cells.Item(151, 4).Value = 45215
$ws.Cells.Item(151, 5).Value = 16
$ws.Cells.Item(151, 6).Value = "Fruta"
$ws.Cells.Item(151, 7).Value = 100108
$ws.Cells.Item(151, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(151, 9).Value = 100108002
$ws.Cells.Item(151, 10).Value = "Mango"
$ws.Cells.Item(151, 11).Value = "Sin especificar"
$ws.Cells.Item(151, 12).Value = "Primera"
$ws.Cells.Item(151, 13).Value = 50
$ws.Cells.Item(151, 14).Value = 10000
$ws.Cells.Item(151, 15).Value = 10000
$ws.Cells.Item(151, 16).Value = 10000
$ws.Cells.Item(151, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(151, 18).Value = "Brasil"
$ws.Cells.Item(151, 19).Value = 2500
$ws.Cells.Item(151, 20).Value = 4
